# Checkpoint of Excel Importer
# Removed java-like mod operator (%)
# Implemented % operator (for Excel)
# Re-triggered unit tests accordingly
#
# Adds new formula-driven sample cells (basic math, power, complex power,
# factorial) on "Sample 2" and registers the defined names that point at
# them, matching the updated sample3.xlsx fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New formulas -----------------------------------------------------
# C6 first, since BasicMathVal (added below) refers to it and C3/C10 use
# that name in their own formulas.
$ws.Range("C6").Formula  = "=(1+2-3*4/6+MOD(5,2))^2"
$ws.Range("C4").Formula  = "=POWER(3,4)"
$ws.Range("C5").Formula  = "=POWER((1+2-3*4/6+MOD(5,2))^2,3)"
$ws.Range("C10").Formula = "=FACT(BasicMathVal)"
$ws.Range("C3").Formula  = "=BasicMathVal*50%"

# --- New defined names --------------------------------------------------
$wb.Names.Add('BasicMathVal',    '=''Sample 2''!$C$6')
$wb.Names.Add('ComplexPowerVal', '=''Sample 2''!$C$5')
$wb.Names.Add('FactorialVal',    '=''Sample 2''!$C$10')
$wb.Names.Add('PercentageVal',   '=''Sample 2''!$C$3')
$wb.Names.Add('PowerVal',        '=''Sample 2''!$C$4')
